$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 82
$ws.Cells.Item($row, 1).Value = "sherry"
$ws.Cells.Item($row, 2).Value = "Grocery"
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0

$ws.Range("B82:G82").Select()
